$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$newFilesTabQuery = @"
SELECT
    f1.file_name AS "File Name",
    s.study_name AS "Study Name",
    s.phs_accession AS "Accession",
    sp.participant_id AS "Participant Id",
    COALESCE((
        SELECT
            REPLACE(GROUP_CONCAT(CASE WHEN rn <= 5 THEN smp.sample_id ELSE NULL END, ', '), ', , ', ', ') ||
            CASE WHEN MAX(rn) > 5 THEN ', ...' ELSE '' END
        FROM (
            SELECT
                smp.sample_id,
                ROW_NUMBER() OVER (ORDER BY smp.sample_id) AS rn
            FROM df_sample smp
            WHERE smp."participant.study_participant_id" = sp.study_participant_id
        ) smp
    ), '') AS "Sample Id",
    f1.file_type AS "File Type",
    gi.library_strategy AS "Library Strategy",
'' AS "Supplementary File"
FROM 
    df_study s
INNER JOIN 
    df_participant sp ON sp."study.phs_accession" = s.phs_accession
INNER JOIN  
    df_sample smp ON smp."participant.study_participant_id" = sp.study_participant_id
INNER JOIN 
    df_file f1 ON f1."sample.sample_id" = smp.sample_id
INNER JOIN
    df_genomic_info gi ON gi."file.file_id" = f1.file_id
INNER JOIN
    df_diagnosis d ON d."participant.study_participant_id" = sp.study_participant_id
INNER JOIN
    df_program p ON p.program_acronym = s."program.program_acronym"
WHERE 
    s.phs_accession = 'phs001437' AND gi.library_source = 'Transcriptomic'
GROUP BY
    f1.file_name,
    s.study_name,
    s.phs_accession,
    sp.participant_id,
    f1.file_type,
    gi.library_strategy
ORDER BY 
    f1.file_name ASC
LIMIT 100;
"@
$ws.Range("B4").Value = $newFilesTabQuery

# Leave the cursor on the cell that was just edited, matching the saved
# selection state of the workbook after the edit.
[void]$ws.Range("B4").Select()
